# Applies the "PO Forecast" sheet addition + header-renames described by
# the diff:
#   1. Rename "Weekly Quantity"!B1 from "Requested quantity" to "Weekly_PO_Qty"
#   2. Rename "Monthly Trend"!B1 from "Requested quantity" to "Monthly_PO_Qty"
#   3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with
#      columns ds / PO_Forecast / yhat_lower / yhat_upper and 28 data rows.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: header renames on the existing sheets -------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: new "PO Forecast" sheet --------------------------------------------
$sheetCount = $wb.Worksheets.Count
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$wsForecast.Name = "PO Forecast"

# Match the page margins used by the sibling sheets (0.75in/0.75in/1in/1in,
# 0.5in header/footer == 54/54/72/72/36/36 points).
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row - reuse the bold/bordered/centered header style already used
# by the other two sheets (cellXf index 1) by copy/paste-special of formats
# from an existing header cell, tiled across the 4 destination columns.
$wsWeekly.Range("A1:B1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows - 28 rows (rows 2-29). Column A needs the same date-number-format
# style (cellXf index 2) used by column A on the other two sheets, so we
# copy/paste-special that formatting down the whole column first.
$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A29").PasteSpecial(-4122) | Out-Null

$data = @(
    @(45081.99999999999, 0, -30.65212395824164, 24.4948957939655),
    @(45088.99999999999, 0, -26.97617308290547, 24.75927173340467),
    @(45151.99999999999, 3, -26.27712384241473, 26.70643441762786),
    @(45179.99999999999, 4, -21.84253491075988, 31.52736384175298),
    @(45193.99999999999, 5, -23.22801348831307, 31.21061004722832),
    @(45200.99999999999, 6, -20.5495155227786, 31.96057824160141),
    @(45207.99999999999, 6, -19.49002691876971, 32.9611805952657),
    @(45228.99999999999, 7, -21.28839395940092, 33.41971559641038),
    @(45235.99999999999, 8, -19.36810684722652, 33.85645581789311),
    @(45242.99999999999, 8, -18.3937980445621, 34.80962035783814),
    @(45249.99999999999, 8, -18.20526296844842, 33.89699261997431),
    @(45333.99999999999, 13, -12.28267601373856, 39.31696044393701),
    @(45347.99999999999, 14, -11.71815686642807, 41.25887175161218),
    @(45361.99999999999, 15, -13.32727388568488, 40.44702950283356),
    @(45389.99999999999, 16, -10.50724232824384, 42.85067841025483),
    @(45396.99999999999, 17, -10.58628905481849, 44.76491152764937),
    @(45403.99999999999, 17, -9.004838241168171, 44.65788074616678),
    @(45410.99999999999, 18, -8.558817684470556, 45.3053674743157),
    @(45424.99999999999, 18, -6.306750696441921, 45.98852088479056),
    @(45445.99999999999, 20, -7.288721879126082, 47.8513510982541),
    @(45452.99999999999, 20, -8.278877222769102, 46.46579202686188),
    @(45459.99999999999, 20, -8.842509948652884, 47.77565812098833),
    @(45466.99999999999, 21, -6.752901443864514, 46.56533503315993),
    @(45473.99999999999, 21, -4.574108641495878, 46.74094840863094),
    @(45480.99999999999, 22, -3.824539320816278, 48.10975149219672),
    @(45487.99999999999, 22, -5.874019670676055, 48.33135233105263),
    @(45494.99999999999, 22, -3.827604905718416, 48.5705023286648),
    @(45501.99999999999, 23, -4.215619214756894, 48.06439805068436)
)

$rowIdx = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($rowIdx, 1).Value = $row[0]
    $wsForecast.Cells.Item($rowIdx, 2).Value = $row[1]
    $wsForecast.Cells.Item($rowIdx, 3).Value = $row[2]
    $wsForecast.Cells.Item($rowIdx, 4).Value = $row[3]
    $rowIdx++
}

$wsForecast.Range("A1").Select() | Out-Null

Write-Output "PO Forecast sheet added; headers renamed."
